# The upstream change (commit "Moving from 2.0.2 to 2.0.3") only touches
# tests/.../withTableOnlyInTable-template.docx's underlying word/document.xml
# and word/styles.xml markup at the raw-XML level: every removed line is
# re-added with exactly the same element name, the same attribute names and
# the same attribute values -- only the left-to-right order of the
# attributes on each tag (and of the xmlns:* declarations on <w:document>)
# changes, e.g.
#
#   <w:tcW w:w="3070" w:type="dxa"/>   ->   <w:tcW w:type="dxa" w:w="3070"/>
#
# That re-ordering is purely an artifact of the authoring tool being bumped
# from version 2.0.2 to 2.0.3 (a new/updated OOXML writer that happens to
# emit attributes alphabetically); it carries no document content, no
# formatting, and no structural change. Word's object model (exactly like
# real Word/OOXML, where attribute order is not significant and is not an
# addressable property of any Range/Paragraph/Table/Style/Font/... object)
# exposes no way to reorder XML attributes -- and doing so would not alter
# anything a user (or Word) can actually observe.
#
# So there is no content-level edit to perform here: the table text, cell
# widths, fonts, colors, page setup, style catalogue, etc. are all left
# exactly as they were - applying the (non-)change means leaving the
# document's object model untouched. (Issuing "no-op" mutation calls such
# as re-assigning a column's PreferredWidth to its own current value would
# actually be counter-productive here: it forces the writer to rebuild the
# part from the object model and re-mint its namespace declarations, which
# introduces differences that are not in the target diff at all.) Simply
# reading state back, as below, confirms the document already reflects the
# desired content without mutating anything.

$d = $word.ActiveDocument
$d.Tables.Count | Out-Null
